# Fruta / hortaliza, semanal
# Insert two new weekly price rows (96 and 97) into the "Pera" sheet,
# pushing the existing rows 96:126 down to 98:128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (rows 96 and 97).
# This shifts the previous rows 96-126 down to 98-128, exactly matching
# the target workbook's row layout (dimension grows from A1:T126 to A1:T128).
$ws.Range("A96:A97").EntireRow.Insert()

# New row 96: Pera, Packham's Triumph, Especial, 2021-09-21 (serial 44460)
$ws.Range("A96").Value = 7
$ws.Range("B96").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C96").Value = 'Ñuble'
$ws.Range("D96").Value = 44460
$ws.Range("E96").Value = 16
$ws.Range("F96").Value = 'Fruta'
$ws.Range("G96").Value = 100104
$ws.Range("H96").Value = 'Frutos de pepita'
$ws.Range("I96").Value = 100104005
$ws.Range("J96").Value = 'Pera'
$ws.Range("K96").Value = "Packham's Triumph"
$ws.Range("L96").Value = 'Especial'
$ws.Range("M96").Value = 60
$ws.Range("N96").Value = 11000
$ws.Range("O96").Value = 11000
$ws.Range("P96").Value = 11000
$ws.Range("Q96").Value = '$/caja 16 kilos empedrada'
$ws.Range("R96").Value = 'Provincia de Curicó'
$ws.Range("S96").Value = 688
$ws.Range("T96").Value = 16

# New row 97: Pera, Packham's Triumph, Primera, 2021-09-21 (serial 44460)
$ws.Range("A97").Value = 7
$ws.Range("B97").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C97").Value = 'Ñuble'
$ws.Range("D97").Value = 44460
$ws.Range("E97").Value = 16
$ws.Range("F97").Value = 'Fruta'
$ws.Range("G97").Value = 100104
$ws.Range("H97").Value = 'Frutos de pepita'
$ws.Range("I97").Value = 100104005
$ws.Range("J97").Value = 'Pera'
$ws.Range("K97").Value = "Packham's Triumph"
$ws.Range("L97").Value = 'Primera'
$ws.Range("M97").Value = 120
$ws.Range("N97").Value = 9000
$ws.Range("O97").Value = 10000
$ws.Range("P97").Value = 9500
$ws.Range("Q97").Value = '$/caja 16 kilos empedrada'
$ws.Range("R97").Value = 'Provincia de Curicó'
$ws.Range("S97").Value = 594
$ws.Range("T97").Value = 16
